# Refresh the NATMI St6gal1-Cd22 LR-pair sheet with the newly recomputed
# TPM-based statistics. The "Target cluster = ECs" rows are no longer part
# of the analysis, so only the four "Target cluster = Resolving-Mac" rows
# survive (renumbered as rows 2-5) with their values updated in place; the
# old rows 6-9 (the second half of the former 8-row table) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (Sending=ECs, Target=Resolving-Mac) ---
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 4.948843666666667
$ws.Range("H2").Value = 14.846531
$ws.Range("I2").Value = 0.2171174124320646
$ws.Range("J2").Value = 0.2171174124320646
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 15.65485666666667
$ws.Range("N2").Value = 46.96456999999999
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 77.47343826740777
$ws.Range("R2").Value = 697.26094440667
$ws.Range("S2").Value = 0.2171174124320646
$ws.Range("T2").Value = 0.2171174124320646

# --- Update row 3 (Sending=FAPs, Target=Resolving-Mac) ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("G3").Value = 3.027316
$ws.Range("H3").Value = 9.081948000000001
$ws.Range("I3").Value = 0.1328154738371249
$ws.Range("J3").Value = 0.1328154738371249
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.65485666666667
$ws.Range("N3").Value = 46.96456999999999
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 47.39219806470667
$ws.Range("R3").Value = 426.52978258236
$ws.Range("S3").Value = 0.1328154738371249
$ws.Range("T3").Value = 0.1328154738371249

# --- Update row 4 (Sending=MuSCs, Target=Resolving-Mac) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 4.429917333333333
$ws.Range("H4").Value = 13.289752
$ws.Range("I4").Value = 0.1943508935591658
$ws.Range("J4").Value = 0.1943508935591658
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 15.65485666666667
$ws.Range("N4").Value = 46.96456999999999
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 69.34972089851554
$ws.Range("R4").Value = 624.1474880866399
$ws.Range("S4").Value = 0.1943508935591658
$ws.Range("T4").Value = 0.1943508935591658

# --- Update row 5 (Sending=Resolving-Mac, Target=Resolving-Mac) ---
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 10.38732133333333
$ws.Range("H5").Value = 31.161964
$ws.Range("I5").Value = 0.4557162201716447
$ws.Range("J5").Value = 0.4557162201716447
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 15.65485666666667
$ws.Range("N5").Value = 46.96456999999999
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 162.6120266239422
$ws.Range("R5").Value = 1463.50823961548
$ws.Range("S5").Value = 0.4557162201716447
$ws.Range("T5").Value = 0.4557162201716447

# --- Remove now-obsolete rows 6-9 ---
$ws.Range("A6:T9").EntireRow.Delete()
